$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - RandomForestRegressor (name unchanged, values updated)
$ws.Range("B3").Value = 0.9997520920247163
$ws.Range("C3").Value = 0.9996861250796784
$ws.Range("D3").Value = 0.9810990003356123

# Row 4 - GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9995798436907105
$ws.Range("C4").Value = 0.9995324110764762
$ws.Range("D4").Value = 0.9967430691575127

# Row 5 - AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9998616035638871
$ws.Range("C5").Value = 0.9998469115744985
$ws.Range("D5").Value = 0.9996428616360841
